$wb = $excel.ActiveWorkbook

# --- Rename sheets (adds "spe"/"s" infix to the "add_*" sheets) ---
# Renaming via the Worksheet.Name property keeps every dependent reference
# (definedNames, formulas, etc.) in sync automatically.
$wb.Worksheets.Item("add_Polizas").Name  = "add_spe_Polizas"
$wb.Worksheets.Item("add_Canales").Name  = "add_spe_Canales"
$wb.Worksheets.Item("add_Amparos").Name  = "add_spe_Amparos"
$wb.Worksheets.Item("add_Serfi").Name    = "add_spe_Serfi"
$wb.Worksheets.Item("add_Atipicos").Name = "add_s_Atipicos"

# --- Switch the active sheet/selection ---
# Previously "Fechas" was the active tab (selection E13); now the renamed
# "add_s_Atipicos" sheet is active with C15 selected instead.
$atipicos = $wb.Worksheets.Item("add_s_Atipicos")
$atipicos.Activate() | Out-Null
$atipicos.Range("C15").Select() | Out-Null
